$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5; this shifts existing rows 5-18 down to 6-19
$ws.Rows.Item(5).Insert()

# Fill in the new row 5 with the weekly record for Espárragos
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 44467
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = 300000000
$ws.Cells.Item(5, 7).Value = "Espárragos"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 50
$ws.Cells.Item(5, 11).Value = 3000
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 13).Value = 3000
$ws.Cells.Item(5, 14).Value = "`$/kilo"
$ws.Cells.Item(5, 15).Value = "Región del Maule"
$ws.Cells.Item(5, 16).Value = 3000
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
